# Applies the "Updated test cases of all" commit to TestData_AdactinHotelApp.xlsx
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. LoginPositiveTest: selection moves to H7
# ------------------------------------------------------------------
$wsLoginPos = $wb.Worksheets.Item("LoginPositiveTest")
$wsLoginPos.Activate()
$wsLoginPos.Range("H7").Select() | Out-Null

# ------------------------------------------------------------------
# 2. LoginNegativeTest: recovery-email address updated, new blank
#    formatted cell in D2, selection moves to D6
# ------------------------------------------------------------------
$wsLoginNeg = $wb.Worksheets.Item("LoginNegativeTest")
$wsLoginNeg.Activate()
$wsLoginNeg.Hyperlinks.Delete() | Out-Null
$wsLoginNeg.Range("C2").Value = "adactin2020@gmail.com"
$wsLoginNeg.Hyperlinks.Add($wsLoginNeg.Range("C2"), "mailto:adactin2020@gmail.com") | Out-Null
$wsLoginNeg.Range("C2").Style = "Hyperlink"
$wsLoginNeg.Range("D2").Style = "Hyperlink"
$wsLoginNeg.Range("D6").Select() | Out-Null

# ------------------------------------------------------------------
# 3. RecoveryEmail: brand new worksheet inserted right after
#    LoginNegativeTest
# ------------------------------------------------------------------
$wsRecovery = $wb.Worksheets.Add($null, $wsLoginNeg)
$wsRecovery.Name = "RecoveryEmail"
$wsRecovery.Range("A1").Value = "Email"
$wsRecovery.Range("A2").Value = "adactin2020@gmail.com"
$wsRecovery.Hyperlinks.Add($wsRecovery.Range("A2"), "mailto:adactin2020@gmail.com") | Out-Null
$wsRecovery.Range("A2").Style = "Hyperlink"
$wsRecovery.Range("A3").Value = "Adtraining@adactin.com"
$wsRecovery.Range("A6").Select() | Out-Null

# ------------------------------------------------------------------
# 4. ChangePassword: test data rewritten, sheet becomes active tab
# ------------------------------------------------------------------
$wsChangePwd = $wb.Worksheets.Item("ChangePassword")
$wsChangePwd.Range("A1:E4").ClearContents() | Out-Null

$wsChangePwd.Range("A1").Value = "UserName"
$wsChangePwd.Range("B1").Value = "Password"
$wsChangePwd.Range("C1").Value = "CurrentPassword"
$wsChangePwd.Range("D1").Value = "NewPassword"
$wsChangePwd.Range("E1").Value = "RePassword"

$wsChangePwd.Range("A2").Value = "AdactinTrainee"
$wsChangePwd.Range("B2").Value = "adactin"
$wsChangePwd.Range("C2").Value = "adactin"
$wsChangePwd.Range("D2").Value = "abcd"
$wsChangePwd.Range("E2").Value = "abcd"

$wsChangePwd.Range("A3").Value = "AdactinTrainee"
$wsChangePwd.Range("B3").Value = "adactin"
$wsChangePwd.Range("C3").Value = "adactin"

$wsChangePwd.Range("A4").Value = "AdactinTrainee"
$wsChangePwd.Range("B4").Value = "abcd"
$wsChangePwd.Range("C4").Value = "abcd"
$wsChangePwd.Range("D4").Value = "adactin"
$wsChangePwd.Range("E4").Value = "adac"

# ------------------------------------------------------------------
# 5. Reorder sheets: RecoveryEmail already sits after LoginNegativeTest;
#    move NewUserRegistration to sit after ChangePassword
# ------------------------------------------------------------------
$wsNewUser = $wb.Worksheets.Item("NewUserRegistration")
$wsNewUser.Move($null, $wsChangePwd)

# re-fetch references: moving a sheet can invalidate older handles
$wsChangePwd = $wb.Worksheets.Item("ChangePassword")
$wsNewUser = $wb.Worksheets.Item("NewUserRegistration")
$wsNewUser.Activate()
$wsNewUser.Range("D11").Select() | Out-Null

# ------------------------------------------------------------------
# 6. SearchHotelNegativeTest: selection moves to C16
# ------------------------------------------------------------------
$wsSearchNeg = $wb.Worksheets.Item("SearchHotelNegativeTest")
$wsSearchNeg.Activate()
$wsSearchNeg.Range("C16").Select() | Out-Null

# ------------------------------------------------------------------
# 7. SearchHotelPositiveTest: selection moves to H10
# ------------------------------------------------------------------
$wsSearchPos = $wb.Worksheets.Item("SearchHotelPositiveTest")
$wsSearchPos.Activate()
$wsSearchPos.Range("H10").Select() | Out-Null

# ------------------------------------------------------------------
# 8. BookHotelPositiveTest: selection unchanged, just re-touch so the
#    worksheet dimension gets recalculated
# ------------------------------------------------------------------
$wsBookPos = $wb.Worksheets.Item("BookHotelPositiveTest")
$wsBookPos.Activate()
$wsBookPos.Range("A2:B2").Select() | Out-Null

# ------------------------------------------------------------------
# 9. BookHotelNegativeTest: ExpYear for last row corrected 2019 -> 2016,
#    selection moves to G12
# ------------------------------------------------------------------
$wsBookNeg = $wb.Worksheets.Item("BookHotelNegativeTest")
$wsBookNeg.Activate()
$wsBookNeg.Range("Q11").Value = 2016
$wsBookNeg.Range("G12").Select() | Out-Null

# ------------------------------------------------------------------
# Final active sheet = ChangePassword (matches workbookView activeTab)
# ------------------------------------------------------------------
$wsChangePwd.Activate()
$wsChangePwd.Range("G6").Select() | Out-Null
